$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.116.73'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.483.83'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.84%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.87'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.78'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.17%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +2.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.483.24'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.75%  '
$ws.Range('E10').Value = '  +5.58%  '
$ws.Range('E11').Value = '  +1.75%  '
$ws.Range('E12').Value = '  +4.02%  '
$ws.Range('E13').Value = '  +3.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.44'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.912.02'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.952.93'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000170'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.555.98'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.04'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.50%  '
$ws.Range('E20').Value = '  +1.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '349.42'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.69%  '
$ws.Range('E22').Value = '  +2.38%  '
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.24'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '68.36'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.56%  '
$ws.Range('E26').Value = '  +4.51%  '
$ws.Range('E27').Value = '  +5.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.621.60'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0909'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '512.68'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.73'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.86%  '
$ws.Range('E33').Value = '  +3.61%  '
$ws.Range('E34').Value = '  +0.91%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '160.51'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.80%  '
$ws.Range('E37').Value = '  +5.00%  '
$ws.Range('E38').Value = '  +0.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.25'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.71%  '
$ws.Range('E40').Value = '  +1.83%  '
$ws.Range('E41').Value = '  +3.61%  '
$ws.Range('E42').Value = '  +0.22%  '
$ws.Range('E43').Value = '  +3.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.81'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.21%  '
$ws.Range('E45').Value = '  +3.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '142.95'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.39%  '
$ws.Range('E48').Value = '  +1.77%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.516'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.84%  '
$ws.Range('E50').Value = '  +4.48%  '
$ws.Range('E51').Value = '  +1.44%  '
